$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A34").Value = "Leetcode"
$ws.Range("B34").Value = 981
$ws.Range("C34").Value = "Time Based Key-Value Store"
$ws.Range("D34").Value = "Hashmaps, Binary Search"
$ws.Range("E34").Value = "Medium"
$ws.Range("F34").Value = "Neetcode 150"
$ws.Range("G34").Value = "SOLVED"
$ws.Range("H34").Value = "17/06/2025"
$ws.Range("I34").Value = "Smashed it."
